$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "iaest-measure:estado-civil"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-measure:edad"

# Row 3 updates
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4 updates
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "URI-comarca"

# Remove row 5 entirely (mapping-*.xlsx references no longer apply)
$ws.Range("A5:I5").Delete()
